$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 4, pushing the current
# rows 4 and 5 down to rows 6 and 7 (formatting carries along with them).
$ws.Range("A4:A5").EntireRow.Insert()

# Row 4: new "Segunda" quality entry for Provincia de Chacabuco
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C4").Value = "Los Lagos"
$ws.Range("D4").Value = 45072
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107011
$ws.Range("J4").Value = "Tuna"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 16000
$ws.Range("Q4").Value = "$/caja 18 kilos"
$ws.Range("R4").Value = "Provincia de Chacabuco"
$ws.Range("S4").Value = 889
$ws.Range("T4").Value = 18

# Row 5: new "Segunda" quality entry for Provincia de Limarí
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C5").Value = "Los Lagos"
$ws.Range("D5").Value = 45072
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = "Otros"
$ws.Range("I5").Value = 100107011
$ws.Range("J5").Value = "Tuna"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 17000
$ws.Range("P5").Value = 17000
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 944
$ws.Range("T5").Value = 18
